$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A251").Value = '17u145t781k4f4h0ytt8rdo03ph1phl9'
$ws.Range("B251").Value = '2018-11-18T13:44:52.198'
$ws.Range("A252").Value = '61k0re8et37823gsdod2kadx07f6183p'
$ws.Range("B252").Value = '2018-11-18T13:45:21.541'
$ws.Range("A253").Value = 'nw1lpd99ordd1hst1w0fmuc98ueo6bs1'
$ws.Range("B253").Value = '2018-11-18T13:46:23.893'
$ws.Range("A254").Value = 'g0ppkhet375tk5n9dr6yzwc06k727u3t'
$ws.Range("B254").Value = '2018-11-18T13:47:19.423'
$ws.Range("A255").Value = '1ck483036u27qy716uk27531porurj2v'
$ws.Range("B255").Value = '2018-11-18T13:48:19.546'
$ws.Range("A256").Value = 'gb711nw27w17f171a3d12tz7g4t17les'
$ws.Range("B256").Value = '2018-11-18T13:49:19.036'
$ws.Range("A257").Value = 'p9o1thjcg5jsevwucr5d92xfmv3oou64'
$ws.Range("B257").Value = '2018-11-18T14:18:17.126'
$ws.Range("A258").Value = 'bv103ve3ib9z64b224ag2j1qqi6s1fh1'
$ws.Range("B258").Value = '2018-11-18T15:28:07.658'
$ws.Range("A259").Value = 'gp9015o381pi2mg0jxhx7cu476c3yp3k'
$ws.Range("B259").Value = '2018-11-18T15:31:35.670'
$ws.Range("A260").Value = 'muyc735buik6h13apf184q0p3k511y5j'
$ws.Range("B260").Value = '2018-11-18T15:32:20.710'
$ws.Range("A261").Value = '8p287960cc7tebo1o02119d2rffct6xq'
$ws.Range("B261").Value = '2018-11-18T15:33:21.771'
$ws.Range("A262").Value = '7s3o2t6w693sgvd5r0071573m0rm0tm0'
$ws.Range("B262").Value = '2018-11-18T15:34:22.426'
$ws.Range("A263").Value = 'icksbpws017y74le6pi462cuslw9km89'
$ws.Range("B263").Value = '2018-11-18T15:35:22.066'
$ws.Range("A264").Value = 'nwdm8tn8x71747ju8nvdb2pv1x3cgyfv'
$ws.Range("B264").Value = '2018-11-18T15:36:22.097'
$ws.Range("A265").Value = '1m9p79q00z6xc30z7037110916401pc6'
$ws.Range("B265").Value = '2018-11-18T15:37:17.320'
$ws.Range("A266").Value = '8099h0tpr7663rrqgbn5cu918bmmhv45'
$ws.Range("B266").Value = '2018-11-18T15:37:30.943'
$ws.Range("A267").Value = '9of0uz0zx9750a58rzywblg8725m5c38'
$ws.Range("B267").Value = '2018-11-18T15:37:46.403'
$ws.Range("A268").Value = 'k6u21iuy23xs5m6hv542dtr9cujbf9r6'
$ws.Range("B268").Value = '2018-11-18T15:38:26.024'
$ws.Range("A269").Value = '83c2rmuhmaeb915wump31ynema827s7s'
$ws.Range("B269").Value = '2018-11-18T15:39:21.642'
$ws.Range("A270").Value = 'a6sxn817g34630o5e51w17i4zqdhxhti'
$ws.Range("B270").Value = '2018-11-18T15:40:24.605'
$ws.Range("A271").Value = '346r60m4w7dvcnwow693u145r92b2e98'
$ws.Range("B271").Value = '2018-11-18T15:41:21.638'
$ws.Range("A272").Value = 'vu0v66i6p8pv523b03psv5295r4ze598'
$ws.Range("B272").Value = '2018-11-18T15:41:35.639'
$ws.Range("A273").Value = '758u354udtp8tm8kytg3awg2ef1l2wd5'
$ws.Range("B273").Value = '2018-11-18T15:41:50.058'
$ws.Range("A274").Value = 'olzrb7m58kzx58xuzy1k5i9f8q5714d5'
$ws.Range("B274").Value = '2018-11-18T15:42:04.346'
$ws.Range("A275").Value = 'bh301ehyxb0m4m9q7130e125fp1p86pr'
$ws.Range("B275").Value = '2018-11-18T15:45:04.989'
$ws.Range("A276").Value = 'b5n3vv7y1am9sk0qyi5eg8plhn2981l6'
$ws.Range("B276").Value = '2018-11-18T15:46:20.346'
$ws.Range("A277").Value = '9a4i6tbafg7f64mw0fe95fprv6k7bvi4'
$ws.Range("B277").Value = '2018-11-18T15:47:21.474'
$ws.Range("A278").Value = '434swwn50350ufy339vrca6w575a5ojs'
$ws.Range("B278").Value = '2018-11-18T15:48:22.572'
$ws.Range("A279").Value = 'l1hhystxc341gj2tp26rjhb7for47he2'
$ws.Range("B279").Value = '2018-11-18T15:49:24.688'
$ws.Range("A280").Value = 'wk5n2c5f5bo1guqqp0942u59srw7n08l'
$ws.Range("B280").Value = '2018-11-18T15:50:19.727'
$ws.Range("A281").Value = 'tzz28fi28stv5h1t58pcmcs79i85k91d'
$ws.Range("B281").Value = '2018-11-18T15:51:19.734'
$ws.Range("A282").Value = '164z658pale6u6chvx9wh1oc9wy68evi'
$ws.Range("B282").Value = '2018-11-18T15:51:32.663'
$ws.Range("A283").Value = 'f4rsh79h5uj3hh58i74ivrvm020fupx3'
$ws.Range("B283").Value = '2018-11-18T15:51:45.771'
$ws.Range("A284").Value = 'jz8t65n972od3vk6t528f7cduf5os62l'
$ws.Range("B284").Value = '2018-11-18T15:52:20.572'
$ws.Range("A285").Value = '4e8n4t09945a806the6470q33h040856'
$ws.Range("B285").Value = '2018-11-18T15:53:21.043'
$ws.Range("A286").Value = '3bvn1zi4ff1gtigse438mt0641ue0pqx'
$ws.Range("B286").Value = '2018-11-18T15:54:20.889'
$ws.Range("A287").Value = 'iy7l018iqcuewvw7htzqphmq0d2nbdb2'
$ws.Range("B287").Value = '2018-11-18T15:54:35.176'
$ws.Range("A288").Value = '4a6j74zzsg3l9670a5hg4x9tas939a2n'
$ws.Range("B288").Value = '2018-11-18T15:54:49.729'
$ws.Range("A289").Value = 'w7di2onaj13jg59m0vyagyls5ie90ep8'
$ws.Range("B289").Value = '2018-11-18T15:55:04.112'
$ws.Range("A290").Value = 'yu9ivh0tzfkp32gi510njw6vducgxat3'
$ws.Range("B290").Value = '2018-11-18T15:55:18.340'
$ws.Range("A291").Value = 's0cgc244hzej8lbqhau08hur3u6q81is'
$ws.Range("B291").Value = '2018-11-18T16:02:35.231'
$ws.Range("A292").Value = 'hi67x26zandpj3w06yu0kb27g3c07057'
$ws.Range("B292").Value = '2018-11-18T16:03:22.739'
$ws.Range("A293").Value = '3t94q6n1ozmd43ly3fl391tyx0zp45ov'
$ws.Range("B293").Value = '2018-11-18T16:04:22.800'
$ws.Range("A294").Value = '0cc4bp76ubqd0y7lq47d3e71ze0t0c0v'
$ws.Range("B294").Value = '2018-11-18T16:05:23.926'
$ws.Range("A295").Value = 'uxjk5j0b13tr2aa9n91s9o7d0xo7p9e2'
$ws.Range("B295").Value = '2018-11-18T16:06:20.929'
$ws.Range("A296").Value = 'yrvtns4ui6v57iac45ri9rkh4vu934df'
$ws.Range("B296").Value = '2018-11-18T16:07:22.010'
$ws.Range("A297").Value = 'x333aez1xm8r3pvs3tcs2xcod0g2azt5'
$ws.Range("B297").Value = '2018-11-18T16:08:23.681'
$ws.Range("A298").Value = '5a589dg78il0i3kd0bze3cm9znwx7o08'
$ws.Range("B298").Value = '2018-11-18T16:08:37.954'
$ws.Range("A299").Value = 'tdlv76kneitdzvm64862738j2b02fo80'
$ws.Range("B299").Value = '2018-11-18T16:08:52.212'
$ws.Range("A300").Value = 'rl8tr2jy131z6fhuo4706dkj95mhn40y'
$ws.Range("B300").Value = '2018-11-18T16:09:22.954'
$ws.Range("A301").Value = 'qpmj84y0614yc3q0viqsb6dg6u1ab64u'
$ws.Range("B301").Value = '2018-11-18T16:10:18.199'
$ws.Range("A302").Value = '91o07d89tq9q4m0t9vq5a6urxi7936bu'
$ws.Range("B302").Value = '2018-11-18T16:10:32.357'
$ws.Range("A303").Value = 'vnl890oox3eqi0pm6bn48y4f1y129400'
$ws.Range("B303").Value = '2018-11-18T16:11:24.729'
$ws.Range("A304").Value = 'y0yzyx1q2j02tml92ex6y6705g1j7t79'
$ws.Range("B304").Value = '2018-11-18T16:11:55.075'
$ws.Range("A305").Value = 'dj2uu5l3sd692g38rndxatvb6xjtf50z'
$ws.Range("B305").Value = '2018-11-18T16:12:22.458'
$ws.Range("A306").Value = 'exzq6dvztg0eudi2p1zrd6m2zo845f72'
$ws.Range("B306").Value = '2018-11-18T16:13:20.658'
$ws.Range("A307").Value = '2k901zc80n00ys5p87dx8890pxogodcn'
$ws.Range("B307").Value = '2018-11-18T16:14:22.307'
$ws.Range("A308").Value = '17v67h37j4r3wy8tc1bb3vuj2g53geqv'
$ws.Range("B308").Value = '2018-11-18T16:14:37.106'
$ws.Range("A309").Value = 'po05wk2b0570ez28wf25srd2bt675fl8'
$ws.Range("B309").Value = '2018-11-18T16:14:51.823'
$ws.Range("A310").Value = 'c802tu5xdix0dqjq0cs7qxix4z016wj1'
$ws.Range("B310").Value = '2018-11-18T16:15:06.690'
$ws.Range("A311").Value = '52yf6xzc64nqv6386lotkkt0jwr91z0h'
$ws.Range("B311").Value = '2018-11-18T16:32:33.982'
$ws.Range("A312").Value = 'ak9toi8uuh4t02b170ap7rh282p3n5g2'
$ws.Range("B312").Value = '2018-11-18T16:33:20.646'
$ws.Range("A313").Value = 'n23ay9drpo68zi6krh822n87ge2ytf4f'
$ws.Range("B313").Value = '2018-11-18T16:34:23.268'
$ws.Range("A314").Value = '23dc3bbd2hiu1kz0o8f68efwai0v2nhu'
$ws.Range("B314").Value = '2018-11-18T16:35:19.363'
$ws.Range("A315").Value = 'rpnrv13rxfg87yhjqd49ilvi76ziywlj'
$ws.Range("B315").Value = '2018-11-18T16:36:14.683'
$ws.Range("A316").Value = '18m8befl84wu6uuk5816c8xdie8gy822'
$ws.Range("B316").Value = '2018-11-18T16:37:26.135'
$ws.Range("A317").Value = 'lary9ua94f15c56g0x7hy7a7573796eq'
$ws.Range("B317").Value = '2018-11-18T16:38:19.683'
$ws.Range("A318").Value = 'y5st0w25khqx1pa8ixp7uw1gmjy81kk3'
$ws.Range("B318").Value = '2018-11-18T16:38:33.770'
$ws.Range("A319").Value = 'b53z3rxlpx39723b0rd7tt4v0g8yz829'
$ws.Range("B319").Value = '2018-11-18T16:38:48.868'
$ws.Range("A320").Value = 'zqny38a93g27rtx2hloqv581v1qs22mi'
$ws.Range("B320").Value = '2018-11-18T16:39:20.596'
$ws.Range("A321").Value = '4qb833ft0491c3056i4sqbl0u95r781l'
$ws.Range("B321").Value = '2018-11-18T16:40:23.414'
$ws.Range("A322").Value = '36sjmwe01ob9m24q2dfsi646y8x1c84m'
$ws.Range("B322").Value = '2018-11-18T16:41:21.065'
$ws.Range("A323").Value = 'c4n8y556akax7xqh2bd38896ayzogt9h'
$ws.Range("B323").Value = '2018-11-18T16:42:25.294'
$ws.Range("A324").Value = '8592f723be79f9dxal85ua9w7w3jni35'
$ws.Range("B324").Value = '2018-11-18T16:43:25.448'
$ws.Range("A325").Value = 'tzanx756wvzmbi5j044s0th8725byfm2'
$ws.Range("B325").Value = '2018-11-18T16:44:23.533'
$ws.Range("A326").Value = 'o2yy3d563xey80391qb0xqn2dsngvbyk'
$ws.Range("B326").Value = '2018-11-18T16:45:20.351'
$ws.Range("A327").Value = '7mriqgv4hg7ns0gzz83p16e6ud21zyrs'
$ws.Range("B327").Value = '2018-11-18T16:46:21.303'
$ws.Range("A328").Value = '51kdnu2q13k34jzi0s6n277xz3iq9891'
$ws.Range("B328").Value = '2018-11-18T16:46:38.040'
$ws.Range("A329").Value = 'm0hn7ve7f21y02dm5x2lbmbms9dbgjv0'
$ws.Range("B329").Value = '2018-11-18T16:46:53.091'
$ws.Range("A330").Value = 'my8fun79j02u78rgbz07l7xw03sq954n'
$ws.Range("B330").Value = '2018-11-18T16:47:08.564'
$ws.Range("A331").Value = 'nk5og5jnwx8hfth8on389axo077lb84w'
$ws.Range("B331").Value = '2018-11-18T16:57:03.971'
$ws.Range("A332").Value = 'ul9176ls41ajvilt54r21o88n9p54p6s'
$ws.Range("B332").Value = '2018-11-18T17:03:27.886'
$ws.Range("A333").Value = '570818mbvhbm73kjmu01qfm8p1nlmg5r'
$ws.Range("B333").Value = '2018-11-18T17:03:42.806'
$ws.Range("A334").Value = 't6eh2114my1pv22lg5bvc329c3gdww82'
$ws.Range("B334").Value = '2018-11-18T17:03:58.183'
$ws.Range("A335").Value = '4qsd9u5ooi9032v5wh6260hd9l83bo77'
$ws.Range("B335").Value = '2018-11-18T17:04:16.619'
$ws.Range("A336").Value = 'o3ku82mui3otkkl32xgep3vuz4eom9pb'
$ws.Range("B336").Value = '2018-11-18T17:04:32.519'
$ws.Range("A337").Value = 'om2z9n0g2w31po4kj32g94v8b9aig4pj'
$ws.Range("B337").Value = '2018-11-18T17:04:49.365'
$ws.Range("A338").Value = '83weialxk16tll24y4jy9bw0f2669wwe'
$ws.Range("B338").Value = '2018-11-18T17:05:05.208'
$ws.Range("A339").Value = 'w97864c0k18017210d70s08fq4luxk95'
$ws.Range("B339").Value = '2018-11-18T17:05:22.504'
$ws.Range("A340").Value = 'vm87d8t22m5sn6e4f80jl2pi62j7hvo9'
$ws.Range("B340").Value = '2018-11-18T17:05:37.939'
$ws.Range("A341").Value = '24zqime6crdwm99dc2gc69zy6tjpuosw'
$ws.Range("B341").Value = '2018-11-18T17:05:54.402'
$ws.Range("A342").Value = 'lg63s7ngi9pjfo732v5h0h6v4i0slmyi'
$ws.Range("B342").Value = '2018-11-18T17:06:10.410'
$ws.Range("A343").Value = 'qdtn77qpv9cfy8hmq7du2k84t6ssmhtj'
$ws.Range("B343").Value = '2018-11-18T17:06:27.449'
$ws.Range("A344").Value = 'w5bfw15xd4258vr2159ll42g382p918a'
$ws.Range("B344").Value = '2018-11-18T17:06:44.850'
$ws.Range("A345").Value = 'x936a1f6xjvda8hdw5j6fd2ft3t153w6'
$ws.Range("B345").Value = '2018-11-18T17:07:04.817'
$ws.Range("A346").Value = '5i66067brcglbm5nfrv0c59s08s35sq0'
$ws.Range("B346").Value = '2018-11-18T17:07:19.979'
$ws.Range("A347").Value = 'ua709uib7z6z69hwd14dkmo1o29js64j'
$ws.Range("B347").Value = '2018-11-18T17:07:39.878'
$ws.Range("A348").Value = '84lhgbwghk8mjb5ivmw4h4va2zzs5n75'
$ws.Range("B348").Value = '2018-11-18T17:07:59.141'
$ws.Range("A349").Value = 'p56f4vv7tm7532va85zl873r3p4adp08'
$ws.Range("B349").Value = '2018-11-18T17:08:19.990'
$ws.Range("A350").Value = '5o5w8e5eo0o8z9a5ux57g2vx17ki3298'
$ws.Range("B350").Value = '2018-11-18T17:08:40.397'
$ws.Range("A351").Value = 'aelc492u53etxfo2oiq8zq5bc1mnkf69'
$ws.Range("B351").Value = '2018-11-18T17:08:58.718'
$ws.Range("A352").Value = 'd5f6d9m392i5th4i8i642a592p8o7e65'
$ws.Range("B352").Value = '2018-11-18T17:14:00.682'
$ws.Range("A353").Value = 'ft320y7g3u85qb8lizu4el90kdxzf4ay'
$ws.Range("B353").Value = '2018-11-18T17:14:21.451'
$ws.Range("A354").Value = '1jf7jsfuvadhoky84fonzq147z58419x'
$ws.Range("B354").Value = '2018-11-18T17:15:24.255'
$ws.Range("A355").Value = 'vwj9tz7d30i9gd73na0taylkmflj40e3'
$ws.Range("B355").Value = '2018-11-18T17:16:20.564'
$ws.Range("A356").Value = '9c9506c6k5klhk58a5w62p7ch34uveje'
$ws.Range("B356").Value = '2018-11-18T17:17:22.364'
$ws.Range("A357").Value = 'tmaf9i38q53e71g3vn3y2j8z4rfp5448'
$ws.Range("B357").Value = '2018-11-18T17:18:27.485'
$ws.Range("A358").Value = 'bd1753tkfh8yn2llgah9fr8yq303grj2'
$ws.Range("B358").Value = '2018-11-18T17:19:23.895'
$ws.Range("A359").Value = 'g9131e6oyl4tb5udb7vku7apr0b5ev0l'
$ws.Range("B359").Value = '2018-11-18T17:19:38.387'
$ws.Range("A360").Value = 'zr4c308d419hhug1lu488a0pp0135u9n'
$ws.Range("B360").Value = '2018-11-18T17:19:52.727'
$ws.Range("A361").Value = 'g2klp2g63h2o0fyust87tnqer962ha98'
$ws.Range("B361").Value = '2018-11-18T17:20:19.223'
$ws.Range("A362").Value = 'alzulzv5x3jkgiqcjjrs38qfn9a2r10f'
$ws.Range("B362").Value = '2018-11-18T17:21:24.948'
$ws.Range("A363").Value = 'ya1e8j516sxux5a159mqvjrq5rsffl4e'
$ws.Range("B363").Value = '2018-11-18T17:22:25.715'
$ws.Range("A364").Value = 'kynt14ns11b8ovvyz5lzcxv7oe8f19m5'
$ws.Range("B364").Value = '2018-11-18T17:24:56.217'
$ws.Range("A365").Value = '8r4o2qwdikbhgjc8g3t524e8ihhg0nt5'
$ws.Range("B365").Value = '2018-11-18T17:25:17.850'
$ws.Range("A366").Value = 'sp0tbe18720cdz41d30sjsmayc9x72pd'
$ws.Range("B366").Value = '2018-11-18T17:26:18.612'
$ws.Range("A367").Value = '57v9um019g316s170orbmc7t5jd028w8'
$ws.Range("B367").Value = '2018-11-18T17:27:23.880'
$ws.Range("A368").Value = '10dp7q74fbm6qmtsu8mawe208yh3m6t8'
$ws.Range("B368").Value = '2018-11-18T17:28:23.629'
$ws.Range("A369").Value = 'tfuyiktzdz33j4u2fz3h9g9ugw1dxm98'
$ws.Range("B369").Value = '2018-11-18T17:28:39.937'
$ws.Range("A370").Value = '95xq8rgxb12j200s5u5py603c2e2460x'
$ws.Range("B370").Value = '2018-11-18T17:28:54.806'
$ws.Range("A371").Value = '6o21hu35gf931munpncl662lql8c9fq5'
$ws.Range("B371").Value = '2018-11-18T17:29:09.787'
$ws.Range("A372").Value = 'jfi8n0r513yk0ib724alt2hhff351699'
$ws.Range("B372").Value = '2018-11-18T17:30:44.332'
$ws.Range("A373").Value = '50h741465o4q0mp3el9f1u713a94x09m'
$ws.Range("B373").Value = '2018-11-18T17:30:58.118'
$ws.Range("A374").Value = 'oe9g01ylr09455e9162w047j7534zvb5'
$ws.Range("B374").Value = '2018-11-18T17:31:26.736'
$ws.Range("A375").Value = '75nngn8ucoqhb3nxtgq25wiay5ths7c7'
$ws.Range("B375").Value = '2018-11-18T17:32:18.941'
